# Append a new data row (row 49) to Sheet1, matching the run performed on 2026-01-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date formatted as text (e.g. "01/11/2026" on the row above).
# Prefixing with an apostrophe forces Excel to store it as literal text instead
# of auto-converting it into a date serial number.
$ws.Range("A49").Value = "'01/12/2026"

$ws.Range("B49").Value = 12844.42
$ws.Range("C49").Value = 0.2147146019633982
$ws.Range("D49").Value = 0.7852853980366018
$ws.Range("E49").Value = -124.01
$ws.Range("F49").Value = -19.53
$ws.Range("G49").Value = -20605.22
$ws.Range("H49").Value = -67.14
$ws.Range("I49").Value = -394.97
$ws.Range("J49").Value = -12.53
